# "Start more entry on airlines"
# Fill in the next two data points in the $INDIGO "Profile" block:
#   row 24 -> HQ      -> "Gurgaon, India"
#   row 25 -> Founded -> 2005
# then leave the selection on the next empty (merged) row, ready for
# the next entry - matching the author's in-progress editing state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets("Main")

$ws.Range("C24").Value = "Gurgaon, India"
$ws.Range("C25").Value = 2005

$ws.Range("C26:D26").Select()
